# Update the "取得日時" (acquisition datetime) timestamps in the "ランサーズ"
# sheet for rows 2-9, column A, from the old scrape time to the new one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2026-01-17 12:48:03"

for ($row = 2; $row -le 9; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
